# Fixed Stimulus Absolute Timestamps
$wb = $excel.ActiveWorkbook

# Rename sheets with updated timestamp-based task order identifiers
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778882490215"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778903320167"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778903330176"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778903966825"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650477890459681"

# Sheet1 (GNG) stimulus file names
$ws1.Range("B2").Value = "go_stims-16504778882090604.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778882320514.csv"
$ws1.Range("B4").Value = "go_stims-16504778882330186.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778882480166.csv"

# Sheet2 (NB) stimulus file names
$ws2.Range("B2").Value = "ZB-match_0-1650477888894017.csv"
$ws2.Range("B3").Value = "TB-1650477889923051.csv"
$ws2.Range("B4").Value = "OB-16504778891080155.csv"
$ws2.Range("B5").Value = "ZB-match_5-16504778885860205.csv"
$ws2.Range("B6").Value = "OB-16504778892040157.csv"
$ws2.Range("B7").Value = "ZB-match_9-16504778885300233.csv"
$ws2.Range("B8").Value = "OB-165047788923802.csv"
$ws2.Range("B9").Value = "TB-16504778903190525.csv"
$ws2.Range("B10").Value = "TB-16504778895550544.csv"

# Sheet4 (TOL) stimulus file names
$ws4.Range("B2").Value = "MM_stims-16504778903640506.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778903400512.csv"
$ws4.Range("B4").Value = "MM_stims-16504778903796813.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477890365016.csv"
$ws4.Range("B6").Value = "MM_stims-16504778903957171.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477890380685.csv"

# Sheet5 (vSAT) stimulus file names
$ws5.Range("B2").Value = "SAT_stims-1650477890398683.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778904276807.csv"
$ws5.Range("B4").Value = "SAT_stims-1650477890411682.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477890443681.csv"
